$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change B2/C2 from "X" text to boolean TRUE, add D2 = FALSE
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = $false

# Row 3: add C3 = "O"
$ws.Range("C3").Value = "O"

# Row 4: add B4 = "O", D4 = "O"
$ws.Range("B4").Value = "O"
$ws.Range("D4").Value = "O"

# Row 5 (new row)
$ws.Range("A5").Value = 45337
$ws.Range("B5").Value = "O"
$ws.Range("C5").Value = "X"
$ws.Range("D5").Value = "O"

# Apply the custom date format (mm/dd/yyyy) to the whole date column (A1:A5),
# matching it to the header cell too, in a single pass so they all share one style.
$ws.Range("A1:A5").NumberFormat = "mm/dd/yyyy"

# Widen column A to fit the new format/header (closest achievable grid value
# to the authored 10.7109375 stored width)
$ws.Columns.Item(1).ColumnWidth = 9.86

# Update selection to B3
$ws.Range("B3").Select()
